$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.884062051773071
$ws.Range("B1").Value = 2.701278924942017
$ws.Range("C1").Value = 5.145395755767822
$ws.Range("D1").Value = 3.654800176620483
$ws.Range("E1").Value = 0.8789551258087158
